$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("disponible") used to hold a shared string "1" for every row;
# replace each with the real numeric stock count. Setting .NumberFormat
# first (inheriting the existing centered/top alignment already used by
# the column) gives the same "#,##0" numeric style Excel created here.
$ws.Range("F2:F15").NumberFormat = "#,##0"

$ws.Range("F2").Value  = 10
$ws.Range("F3").Value  = 4
$ws.Range("F4").Value  = 7
$ws.Range("F5").Value  = 20
$ws.Range("F6").Value  = 17
$ws.Range("F7").Value  = 15
$ws.Range("F8").Value  = 3
$ws.Range("F9").Value  = 7
$ws.Range("F10").Value = 12
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0

# D12/D13 pick up the same green highlight already used on D14/D15.
$ws.Range("D12:D13").Interior.Color = 5296274

# Five new (empty) rows below the data, formatted like the rest of column F.
$ws.Range("F16:F20").NumberFormat = "#,##0.000"
